# "Improved Data Driven Approach"
#
# The sheet tabs are being relabeled (a cyclic rename), keeping each sheet's
# underlying data/position intact, and the active tab / selection moves from
# the old "Sheet4" tab to the tab that will become the new "Sheet4":
#
#   old name -> new name
#   Sheet2   -> Sheet4
#   Sheet3   -> Sheet1
#   Sheet4   -> Sheet2
#   Sheet1   -> Sheet3
#
# Use unique temporary names first so the cyclic rename never collides with
# an existing sheet name.

$wb = $excel.ActiveWorkbook

$wb.Worksheets.Item("Sheet2").Name = "__tmp_Sheet2__"
$wb.Worksheets.Item("Sheet3").Name = "__tmp_Sheet3__"
$wb.Worksheets.Item("Sheet4").Name = "__tmp_Sheet4__"
$wb.Worksheets.Item("Sheet1").Name = "__tmp_Sheet1__"

$wb.Worksheets.Item("__tmp_Sheet2__").Name = "Sheet4"
$wb.Worksheets.Item("__tmp_Sheet3__").Name = "Sheet1"
$wb.Worksheets.Item("__tmp_Sheet4__").Name = "Sheet2"
$wb.Worksheets.Item("__tmp_Sheet1__").Name = "Sheet3"

# Make the (newly named) "Sheet4" tab the active sheet/tab, with its
# selection moved from F15 to D17.
$activeSheet = $wb.Worksheets.Item("Sheet4")
$activeSheet.Activate()
$activeSheet.Range("D17").Select()

# The (newly named) "Sheet2" tab is no longer the active/selected tab
# (it used to be, as the old "Sheet4"); its own selection (G7) is unchanged.
